$d = $word.ActiveDocument

# Locate the paragraph that ends with "після повної поставки обладнання."
# and insert a new, empty paragraph right after it (before "Умови оплати:").
$target = $d.Content
$target.Find.Execute("після повної поставки обладнання.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$range = $target.Paragraphs.Last.Range
$newRange = $range.InsertParagraphAfter()

# Style the new empty paragraph like the one shown in the diff:
#   <w:spacing w:after="0"/>
#   <w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="uk-UA"/></w:rPr>
$newPara = $d.Paragraphs.Last
$newPara.SpaceAfter = 0
$newPara.Range.Font.Bold = $true
$newPara.Range.Font.Size = 12
$newPara.Range.Font.NameFarEast = "minorHAnsi"
$newPara.Range.LanguageID = "uk-UA"
